# Add a new "Groups" column (L) to the roster sheet and populate it
# with each member's group/schedule info, matching the "Add groups
# reading to excel" change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column L
$ws.Range("L1").Value = "Groups"

# Group/schedule values for the existing test rows.
# Rows 4 and 5 have no group assigned yet, so column L stays blank there.
$ws.Range("L2").Value = "Monday 10:00"
$ws.Range("L3").Value = "Sunday 11:00, Friday 19:00"

# Leave the selection on the last cell that was edited (L3), matching
# where the user's cursor ended up after entering the new data.
$ws.Range("L3").Select()
